# Update data import for NEDC
# Adds a "rem3" helper column (I) that mirrors the existing "rem2" (F) helper
# column, and back-fills the previously-blank F cells (alcgp "120+" rows,
# where rem2 had no numeric code) with the literal "n/a" -- matching what I
# gets for those same rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows (2..89) where column F ("rem2") is currently blank -- these are the
# alcgp = "120+" groups that never got a numeric rem2 code.
$blankRows = @(13,14,15,16,29,30,31,44,45,46,47,60,61,62,63,75,76,77,78,88,89)

# 1) Back-fill the blank F cells with "n/a" first, so the new shared string
#    "n/a" is registered before "rem3" (matching the target string order).
foreach ($r in $blankRows) {
    $ws.Cells.Item($r, 6).Value = "n/a"
}

# 2) New header for column I.
$ws.Cells.Item(1, 9).Value = "rem3"

# 3) Column I mirrors column F for every data row: same numeric code, or
#    "n/a" where F was (and now is) the "n/a" placeholder.
for ($r = 2; $r -le 89; $r++) {
    $fVal = $ws.Cells.Item($r, 6).Value2
    $ws.Cells.Item($r, 9).Value = $fVal
}

# 4) Selection reflects the newly added column I being picked (whole column).
$ws.Range("I1:I1048576").Select()

# 5) Record a sort state: the data (A2:H90) was (re-)sorted ascending by
#    column A, which is a no-op on values but leaves behind the sortState
#    metadata that Excel writes after using Data > Sort.
$sortObj = $ws.Sort
$sortObj.SortFields.Clear()
$sortObj.SortFields.Add($ws.Range("A1:A90"))
$sortObj.SetRange($ws.Range("A2:H90"))
$sortObj.Header = 2
$sortObj.Apply()
